$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data refresh: update Price (D) and Volume(1h) (E) columns for the symbol
# list rows whose values changed in this run. Values are written with a
# leading apostrophe so Excel stores them as literal text (matching the
# original cell type) instead of auto-converting numeric-looking strings
# into Number/Percentage cells.

$ws.Range("D2").Value = "'278.69"
$ws.Range("E2").Value = "'6.81%"
$ws.Range("D3").Value = "'27.29"
$ws.Range("E3").Value = "'1.14%"
$ws.Range("D4").Value = "'4.825"
$ws.Range("E4").Value = "'2.80%"
$ws.Range("D5").Value = "'0.06267"
$ws.Range("E5").Value = "'0.72%"
$ws.Range("D6").Value = "'6.858"
$ws.Range("E6").Value = "'1.62%"
$ws.Range("D7").Value = "'0.8781"
$ws.Range("E7").Value = "'2.82%"
$ws.Range("D8").Value = "'0.9411"
$ws.Range("E8").Value = "'3.02%"
$ws.Range("D9").Value = "'0.1447"
$ws.Range("E9").Value = "'2.99%"
$ws.Range("D10").Value = "'0.05144"
$ws.Range("E10").Value = "'6.24%"
$ws.Range("D11").Value = "'0.07287"
$ws.Range("E11").Value = "'2.68%"
$ws.Range("D12").Value = "'0.03166"
$ws.Range("E12").Value = "'1.92%"
$ws.Range("D13").Value = "'0.09047"
$ws.Range("E13").Value = "'-0.15%"
$ws.Range("D14").Value = "'0.001555"
$ws.Range("E14").Value = "'1.68%"
$ws.Range("D15").Value = "'0.0006282"
$ws.Range("E15").Value = "'1.86%"
$ws.Range("D16").Value = "'0.005902"
$ws.Range("E16").Value = "'-2.04%"
$ws.Range("D17").Value = "'3.449"
$ws.Range("E17").Value = "'0.22%"
$ws.Range("D18").Value = "'3.267"
$ws.Range("E18").Value = "'2.85%"
$ws.Range("E19").Value = "'4.59%"
$ws.Range("D22").Value = "'3.856"
$ws.Range("E22").Value = "'-5.59%"
$ws.Range("D23").Value = "'0.04318"
$ws.Range("E23").Value = "'1.78%"
$ws.Range("D24").Value = "'0.001176"
$ws.Range("E24").Value = "'-3.11%"
$ws.Range("E25").Value = "'4.78%"
$ws.Range("D26").Value = "'0.0001198"
$ws.Range("E26").Value = "'-0.22%"
$ws.Range("D27").Value = "'0.0001688"
$ws.Range("E27").Value = "'2.94%"
$ws.Range("D40").Value = "'0.04029"
$ws.Range("E40").Value = "'2.48%"
$ws.Range("D41").Value = "'0.006414"
$ws.Range("E41").Value = "'55.64%"
$ws.Range("E42").Value = "'3.78%"
$ws.Range("D43").Value = "'0.002102"
$ws.Range("E43").Value = "'-4.95%"
$ws.Range("D44").Value = "'0.01386"
$ws.Range("E44").Value = "'-0.13%"
$ws.Range("D45").Value = "'0.00005140"
$ws.Range("E45").Value = "'-0.45%"
$ws.Range("E46").Value = "'-0.21%"
$ws.Range("D47").Value = "'2.339"
$ws.Range("E47").Value = "'1,151.76%"
$ws.Range("D49").Value = "'0.00002097"
$ws.Range("E49").Value = "'-0.21%"
$ws.Range("D50").Value = "'0.0001997"
$ws.Range("E50").Value = "'-0.21%"
